$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: Jill Hrinda Patten bullet - replace the "CPIR (Center for Parent
# Information & Resources)." tail with
# "Mission Empower Community Parent Resource Center, Executive Director."
# The new text is split across four runs to match the target formatting:
#   - "Mission Empower Community Parent Resource Center" -> default run (no rPr)
#   - ", " -> Public Sans / 20 / 20 (same as surrounding runs)
#   - "Executive Director" -> default run (no rPr)
#   - "." -> default run (no rPr)
# ---------------------------------------------------------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("CPIR (Center for Parent Information & Resources).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $target1 = $d.Range($rng1.Start, $rng1.End)
    $frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                  '<w:p>' +
                    '<w:r><w:t>Mission Empower Community Parent Resource Center</w:t></w:r>' +
                    '<w:r><w:rPr><w:rFonts w:ascii="Public Sans" w:hAnsi="Public Sans"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
                    '<w:r><w:t>Executive Director</w:t></w:r>' +
                    '<w:r><w:t>.</w:t></w:r>' +
                  '</w:p>' +
                '</w:body>' +
              '</w:document>' +
            '</pkg:xmlData>' +
          '</pkg:part>' +
        '</pkg:package>'
    $target1.InsertXML($frag1)
}

# ---------------------------------------------------------------------------
# Edit 2: Barbara Simpson bullet - append a "." run right after the existing
# " Part C Data Manager" run (same Public Sans / 20 / 20 formatting),
# so the sentence now ends with a period.
# ---------------------------------------------------------------------------

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" Part C Data Manager", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $target2 = $d.Range($rng2.Start, $rng2.End)
    $frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                  '<w:p>' +
                    '<w:r><w:rPr><w:rFonts w:ascii="Public Sans" w:hAnsi="Public Sans"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> Part C Data Manager</w:t></w:r>' +
                    '<w:r><w:rPr><w:rFonts w:ascii="Public Sans" w:hAnsi="Public Sans"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r>' +
                  '</w:p>' +
                '</w:body>' +
              '</w:document>' +
            '</pkg:xmlData>' +
          '</pkg:part>' +
        '</pkg:package>'
    $target2.InsertXML($frag2)
}
